$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking price cells as Text so Excel keeps them
# as strings (matching the original inline-string / shared-string cell type)
# instead of silently converting them to numbers.
$textCells = @("D5", "D6", "D12", "D18", "D21", "D22", "D24", "D29", "D31", "D32", "D35", "D38", "D46", "D47", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values from the crypto price/volume refresh.
$ws.Range('D2').Value = '66.596.24'
$ws.Range('D3').Value = '3.071.66'
$ws.Range('E3').Value = '  -1.32%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '574.38'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '169.07'
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.068.69'
$ws.Range('E8').Value = '  -1.30%  '
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('D12').Value = '0.468'
$ws.Range('E12').Value = '  -3.19%  '
$ws.Range('E13').Value = '  -2.75%  '
$ws.Range('E14').Value = '  -3.97%  '
$ws.Range('D16').Value = '3.581.89'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('D17').Value = '66.544.25'
$ws.Range('D18').Value = '16.87'
$ws.Range('E18').Value = '  +3.76%  '
$ws.Range('E19').Value = '  -3.29%  '
$ws.Range('D20').Value = '3.082.04'
$ws.Range('E20').Value = '  -0.96%  '
$ws.Range('D21').Value = '486.25'
$ws.Range('E21').Value = '  +2.13%  '
$ws.Range('D22').Value = '7.69'
$ws.Range('E22').Value = '  -2.71%  '
$ws.Range('E23').Value = '  -3.93%  '
$ws.Range('D24').Value = '82.54'
$ws.Range('E24').Value = '  -1.88%  '
$ws.Range('E25').Value = '  -4.90%  '
$ws.Range('E26').Value = '  -3.49%  '
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = '7.77'
$ws.Range('E29').Value = '  -1.34%  '
$ws.Range('E30').Value = '  -4.90%  '
$ws.Range('D31').Value = '2.59'
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('D32').Value = '27.47'
$ws.Range('E32').Value = '  -3.82%  '
$ws.Range('E33').Value = '  -3.76%  '
$ws.Range('D34').Value = '0.0₃0906'
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('E37').Value = '  -4.74%  '
$ws.Range('D38').Value = '47.07'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('E40').Value = '  -4.92%  '
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('E42').Value = '  -4.84%  '
$ws.Range('D43').Value = '2.764.97'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('E44').Value = '  -2.38%  '
$ws.Range('E45').Value = '  -3.25%  '
$ws.Range('D46').Value = '134.93'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '365.47'
$ws.Range('E47').Value = '  -5.78%  '
$ws.Range('D49').Value = '24.34'
$ws.Range('E49').Value = '  -1.58%  '
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('E51').Value = '  -2.05%  '

# Restore default (Normal) styling on the cells we temporarily reformatted,
# so only the cell values changed, not their appearance/style.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}

